# Updates symbol list / price+volume data for cryptos.xlsx (Sheet1)
# Each cell is rewritten as literal text (leading apostrophe forces Excel to
# store a TEXT value rather than coercing numeric-looking strings like
# "245.84" or "-0.43%" into actual numbers), then the style is reset to
# "Normal" so the quote-prefix formatting introduced by the text-entry does
# not linger as a visible cell-style change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '245.84' }
    @{ Cell = "E2"; Value = '-0.43%' }
    @{ Cell = "D3"; Value = '30.31' }
    @{ Cell = "E3"; Value = '0.71%' }
    @{ Cell = "D4"; Value = '5.155' }
    @{ Cell = "E4"; Value = '-0.46%' }
    @{ Cell = "D5"; Value = '0.05764' }
    @{ Cell = "E5"; Value = '0.62%' }
    @{ Cell = "D6"; Value = '6.667' }
    @{ Cell = "E6"; Value = '1.06%' }
    @{ Cell = "D7"; Value = '3.250' }
    @{ Cell = "E7"; Value = '5.78%' }
    @{ Cell = "D8"; Value = '0.8515' }
    @{ Cell = "E8"; Value = '-1.05%' }
    @{ Cell = "D9"; Value = '0.8571' }
    @{ Cell = "E9"; Value = '-2.65%' }
    @{ Cell = "D10"; Value = '0.1381' }
    @{ Cell = "E10"; Value = '1.09%' }
    @{ Cell = "D11"; Value = '0.07083' }
    @{ Cell = "E11"; Value = '-0.18%' }
    @{ Cell = "D12"; Value = '0.03262' }
    @{ Cell = "E12"; Value = '14.13%' }
    @{ Cell = "D13"; Value = '0.09364' }
    @{ Cell = "D14"; Value = '0.001539' }
    @{ Cell = "E14"; Value = '1.39%' }
    @{ Cell = "B15"; Value = 'TigerCash' }
    @{ Cell = "C15"; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' }
    @{ Cell = "D15"; Value = '0.006147' }
    @{ Cell = "E15"; Value = '2.11%' }
    @{ Cell = "B16"; Value = 'LEO' }
    @{ Cell = "C16"; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' }
    @{ Cell = "D16"; Value = '3.526' }
    @{ Cell = "E16"; Value = '0.80%' }
    @{ Cell = "B17"; Value = 'BTSEToken' }
    @{ Cell = "C17"; Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse' }
    @{ Cell = "D17"; Value = '2.217' }
    @{ Cell = "E17"; Value = '-2.75%' }
    @{ Cell = "B18"; Value = 'One' }
    @{ Cell = "C18"; Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one' }
    @{ Cell = "D18"; Value = '0.01024' }
    @{ Cell = "E18"; Value = '-0.35%' }
    @{ Cell = "D19"; Value = '0.3125' }
    @{ Cell = "E19"; Value = '-1.87%' }
    @{ Cell = "D20"; Value = '0.03393' }
    @{ Cell = "E20"; Value = '3.54%' }
    @{ Cell = "E21"; Value = '1.24%' }
    @{ Cell = "D22"; Value = '3.486' }
    @{ Cell = "E22"; Value = '-0.67%' }
    @{ Cell = "D23"; Value = '0.04136' }
    @{ Cell = "E23"; Value = '-0.22%' }
    @{ Cell = "D24"; Value = '0.1409' }
    @{ Cell = "E24"; Value = '2.19%' }
    @{ Cell = "E25"; Value = '0.96%' }
    @{ Cell = "D26"; Value = '0.004151' }
    @{ Cell = "E26"; Value = '-7.75%' }
    @{ Cell = "D27"; Value = '0.0001199' }
    @{ Cell = "E27"; Value = '-0.84%' }
    @{ Cell = "E28"; Value = '4.50%' }
    @{ Cell = "D40"; Value = '0.03759' }
    @{ Cell = "E40"; Value = '-0.63%' }
    @{ Cell = "D41"; Value = '0.1070' }
    @{ Cell = "E41"; Value = '-0.26%' }
    @{ Cell = "D42"; Value = '0.002199' }
    @{ Cell = "E42"; Value = '-0.02%' }
    @{ Cell = "E43"; Value = '-48.62%' }
    @{ Cell = "D44"; Value = '0.008919' }
    @{ Cell = "E44"; Value = '-11.30%' }
    @{ Cell = "D45"; Value = '0.00005475' }
    @{ Cell = "E45"; Value = '7.34%' }
    @{ Cell = "E46"; Value = '-0.02%' }
    @{ Cell = "D47"; Value = '0.07097' }
    @{ Cell = "E47"; Value = '-20.24%' }
    @{ Cell = "E48"; Value = '-10.85%' }
    @{ Cell = "E49"; Value = '-0.02%' }
    @{ Cell = "E50"; Value = '-0.02%' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $range.Value = "'" + $u.Value
    $range.Style = "Normal"
}

